$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date) for rows 2 through 9: 45224 -> 45233
for ($row = 2; $row -le 9; $row++) {
    $ws.Range("C$row").Value = 45233
}
